$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction12")

# Update the two remaining values
$ws.Range("A1").Value = 24
$ws.Range("B1").Value = 25

# Clear out the rest of the previously populated row (C1:Q1) entirely
$ws.Range("C1:Q1").Clear()
